$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (222nm Basal vs 30 J/m^2-Control)
$ws.Range("C2").Value = 0.390909090909091
$ws.Range("D2").Value = -1.33267531145035
$ws.Range("E2").Value = 2.11449349326853
$ws.Range("F2").Value = 0.943287606242385
$ws.Range("G2").Value = 0.99996960525172

# Row 3 (222nm Basal vs 300 J/m^2-Control)
$ws.Range("C3").Value = 0.254545454545455
$ws.Range("D3").Value = -1.46903894781399
$ws.Range("E3").Value = 1.9781298569049
$ws.Range("F3").Value = 0.987571732951626
$ws.Range("G3").Value = 0.99996960525172

# Row 4 (222nm Basal vs 1000 J/m^2-Control)
$ws.Range("C4").Value = 0.0545454545454545
$ws.Range("D4").Value = -1.66903894781399
$ws.Range("E4").Value = 1.7781298569049
$ws.Range("F4").Value = 0.99996960525172
$ws.Range("G4").Value = 0.99996960525172

# Row 5 (222nm Basal vs 2000 J/m^2-Control)
$ws.Range("C5").Value = 1.41818181818182
$ws.Range("D5").Value = -0.305402584177623
$ws.Range("E5").Value = 3.14176622054126
$ws.Range("F5").Value = 0.134397398839095
$ws.Range("G5").Value = 0.30719405448936

# Row 6 (222nm Suprabasal vs 30 J/m^2-Control)
$ws.Range("F6").Value = 0.928618001972295
$ws.Range("G6").Value = 0.99996960525172

# Row 7 (222nm Suprabasal vs 300 J/m^2-Control)
$ws.Range("F7").Value = 0.986643604872859
$ws.Range("G7").Value = 0.99996960525172

# Row 8 (222nm Suprabasal vs 1000 J/m^2-Control)
$ws.Range("F8").Value = 0.2965222266076
$ws.Range("G8").Value = 0.593044453215201

# Row 9 (222nm Suprabasal vs 2000 J/m^2-Control)
$ws.Range("F9").Value = 0.0176058844871828
$ws.Range("G9").Value = 0.0469490252991542

# Row 10 (254nm Basal vs 30 J/m^2-Control)
$ws.Range("F10").Value = 0.996181161513232
$ws.Range("G10").Value = 0.99996960525172

# Row 11 (254nm Basal vs 300 J/m^2-Control)
$ws.Range("F11").Value = 0.805279658652863
$ws.Range("G11").Value = 0.99996960525172

# Row 12 (254nm Basal vs 1000 J/m^2-Control)
$ws.Range("F12").Value = 0.00276156704646191
$ws.Range("G12").Value = 0.00883701454867811

# Row 13 (254nm Basal vs 2000 J/m^2-Control)
$ws.Range("F13").Value = 0.00000451606832641005
$ws.Range("G13").Value = 0.0000240856977408536

# Row 14 (254nm Suprabasal vs 30 J/m^2-Control)
$ws.Range("F14").Value = 0.481124324802707
$ws.Range("G14").Value = 0.85533213298259

# Row 15 (254nm Suprabasal vs 300 J/m^2-Control)
$ws.Range("F15").Value = 0.000246981671981783
$ws.Range("G15").Value = 0.000987926687927132

# Row 16 (254nm Suprabasal vs 1000 J/m^2-Control)
$ws.Range("F16").Value = 0.0000000000982907089053242
$ws.Range("G16").Value = 0.00000000141969280775811

# Row 17 (254nm Suprabasal vs 2000 J/m^2-Control)
$ws.Range("F17").Value = 0.000000000177461600969764
$ws.Range("G17").Value = 0.00000000141969280775811
